# Append one new row (row 6) to Sheet1 with the ride-log entry captured
# at 2025-05-01T14:17:33.783Z.
#
# Every column in this sheet stores values as literal TEXT, even when the
# text happens to look numeric (e.g. the quantity "2323"), and the notes
# column (A) is sometimes an explicit empty string rather than a truly
# blank cell. Plain `.Value = "2323"` would be auto-coerced by Excel into
# the number 2323, and `.Value = ""` clears the cell outright instead of
# leaving behind an explicit empty string - neither matches the source
# data. For those two cells we instead assign a `="..."` formula, which
# forces Excel to store the result as text (including an explicit empty
# string) without disturbing the cell's format/style. The remaining,
# unambiguous text columns are written with a plain `.Value` assignment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 6

$ws.Range("A$row").Formula = "="""""
$ws.Range("B$row").Value = "أحمد شريم"
$ws.Range("C$row").Formula = "=""2323"""
$ws.Range("D$row").Value = "ايتا"
$ws.Range("E$row").Value = "الرحلة 2"
$ws.Range("F$row").Value = "C2"
$ws.Range("G$row").Value = "NRC"
$ws.Range("H$row").Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٥:١٧:٣٣ م"
